$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "领益智造"
$ws.Range("C8").Value = "成飞集成"
$ws.Range("A9").Value = "天融信"
$ws.Range("A10").Value = "合力泰"
$ws.Range("C10").Value = "步步高"
$ws.Range("A11").Value = "岩山科技"
$ws.Range("C11").Value = "鸿博股份"
$ws.Range("C12").Value = "工业富联"
$ws.Range("C13").Value = "大位科技"
$ws.Range("C14").Value = "吉视传媒"
$ws.Range("C15").Value = "大元泵业"
$ws.Range("C16").Value = "天融信"
$ws.Range("C17").Value = "新易盛"
$ws.Range("C18").Value = "合力泰"
$ws.Range("A19").Value = "奋达科技"
$ws.Range("C19").Value = "指南针"
$ws.Range("A20").Value = "乐鑫科技"
$ws.Range("C20").Value = "四川长虹"
$ws.Range("A21").Value = "中国卫星"
$ws.Range("C21").Value = "金力永磁"
